# "avancée sur le mémoire"
#
# 1) The "Date" placeholder (a fixed/cached datetimeFigureOut field) that lives
#    on the Slide Master and on every Slide Layout is bumped from 19/06/2020
#    to 11/08/2020.
# 2) The "inference" textbox on the only slide gets re-capitalised to
#    "Inference".

$p = $ppt.ActivePresentation

$oldDate = "19/06/2020"
$newDate = "11/08/2020"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        $isDatePlaceholder = $false
        if ($shape.Type -eq 14) {
            try {
                if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePlaceholder = $true
                }
            } catch {
                $isDatePlaceholder = $false
            }
        }
        if ($isDatePlaceholder -and $shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide Master.
$master = $p.Slides.Item(1).Master
Update-DatePlaceholder $master.Shapes

# Every Slide Layout hanging off the master.
for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    Update-DatePlaceholder $layout.Shapes
}

# Slide content: "inference" -> "Inference".
$slide = $p.Slides.Item(1)
$shapes = $slide.Shapes
for ($i = 1; $i -le $shapes.Count; $i++) {
    $shape = $shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "inference") {
            $tr.Text = "Inference"
        }
    }
}
